$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.905.94"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "1.745.36"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5182"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2803"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06119"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").Value = "1.752.85"
$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07041"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6409"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.521"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "25.882.18"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006594"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.23%  "

$ws.Range("D22").Value = "1.973.08"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.141"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.651"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.140"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.516"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.810"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08239"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.668"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.430"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04490"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.612"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9881"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6133"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.93%  "

$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01591"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.920"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3841"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.035"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7241"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05443"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.280"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1122"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.655"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "29.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.81%  "
